$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying data table (rows 2-72) shifts down by one "week" of history:
# each row 48-72 now holds the values previously found in the row above it,
# row 47 gets a new date/volume, and a new row 73 is appended holding the
# data that used to be the last row (72).

# Row 47
$ws.Cells.Item(47, 4).Value = 45072
$ws.Cells.Item(47, 9).Value = "Primera"
$ws.Cells.Item(47, 10).Value = 200
$ws.Cells.Item(47, 11).Value = 18000
$ws.Cells.Item(47, 12).Value = 20000
$ws.Cells.Item(47, 13).Value = 18800
$ws.Cells.Item(47, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(47, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(47, 16).Value = 1253
$ws.Cells.Item(47, 17).Value = 15

# Row 48
$ws.Cells.Item(48, 4).Value = 45062
$ws.Cells.Item(48, 9).Value = "Primera"
$ws.Cells.Item(48, 10).Value = 300
$ws.Cells.Item(48, 11).Value = 18000
$ws.Cells.Item(48, 12).Value = 20000
$ws.Cells.Item(48, 13).Value = 18800
$ws.Cells.Item(48, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(48, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(48, 16).Value = 1253
$ws.Cells.Item(48, 17).Value = 15

# Row 49
$ws.Cells.Item(49, 4).Value = 44754
$ws.Cells.Item(49, 9).Value = "Primera"
$ws.Cells.Item(49, 10).Value = 300
$ws.Cells.Item(49, 11).Value = 17000
$ws.Cells.Item(49, 12).Value = 19000
$ws.Cells.Item(49, 13).Value = 18133
$ws.Cells.Item(49, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(49, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(49, 16).Value = 1209
$ws.Cells.Item(49, 17).Value = 15

# Row 50
$ws.Cells.Item(50, 4).Value = 44790
$ws.Cells.Item(50, 9).Value = "Primera"
$ws.Cells.Item(50, 10).Value = 500
$ws.Cells.Item(50, 11).Value = 15000
$ws.Cells.Item(50, 12).Value = 16000
$ws.Cells.Item(50, 13).Value = 15540
$ws.Cells.Item(50, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(50, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(50, 16).Value = 1036
$ws.Cells.Item(50, 17).Value = 15

# Row 51
$ws.Cells.Item(51, 4).Value = 44819
$ws.Cells.Item(51, 9).Value = "Primera"
$ws.Cells.Item(51, 10).Value = 300
$ws.Cells.Item(51, 11).Value = 15000
$ws.Cells.Item(51, 12).Value = 17000
$ws.Cells.Item(51, 13).Value = 16200
$ws.Cells.Item(51, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(51, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(51, 16).Value = 1080
$ws.Cells.Item(51, 17).Value = 15

# Row 52
$ws.Cells.Item(52, 4).Value = 44400
$ws.Cells.Item(52, 9).Value = "Primera"
$ws.Cells.Item(52, 10).Value = 130
$ws.Cells.Item(52, 11).Value = 24000
$ws.Cells.Item(52, 12).Value = 24000
$ws.Cells.Item(52, 13).Value = 24000
$ws.Cells.Item(52, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(52, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(52, 16).Value = 1600
$ws.Cells.Item(52, 17).Value = 15

# Row 53
$ws.Cells.Item(53, 4).Value = 44791
$ws.Cells.Item(53, 9).Value = "Primera"
$ws.Cells.Item(53, 10).Value = 300
$ws.Cells.Item(53, 11).Value = 16000
$ws.Cells.Item(53, 12).Value = 18000
$ws.Cells.Item(53, 13).Value = 17133
$ws.Cells.Item(53, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(53, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(53, 16).Value = 1142
$ws.Cells.Item(53, 17).Value = 15

# Row 54
$ws.Cells.Item(54, 4).Value = 44810
$ws.Cells.Item(54, 9).Value = "Primera"
$ws.Cells.Item(54, 10).Value = 400
$ws.Cells.Item(54, 11).Value = 17000
$ws.Cells.Item(54, 12).Value = 19000
$ws.Cells.Item(54, 13).Value = 17850
$ws.Cells.Item(54, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(54, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(54, 16).Value = 1190
$ws.Cells.Item(54, 17).Value = 15

# Row 55
$ws.Cells.Item(55, 4).Value = 44449
$ws.Cells.Item(55, 9).Value = "Primera"
$ws.Cells.Item(55, 10).Value = 220
$ws.Cells.Item(55, 11).Value = 22000
$ws.Cells.Item(55, 12).Value = 24000
$ws.Cells.Item(55, 13).Value = 23091
$ws.Cells.Item(55, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(55, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(55, 16).Value = 1539
$ws.Cells.Item(55, 17).Value = 15

# Row 56
$ws.Cells.Item(56, 4).Value = 44399
$ws.Cells.Item(56, 9).Value = "Primera"
$ws.Cells.Item(56, 10).Value = 150
$ws.Cells.Item(56, 11).Value = 22000
$ws.Cells.Item(56, 12).Value = 22000
$ws.Cells.Item(56, 13).Value = 22000
$ws.Cells.Item(56, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(56, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(56, 16).Value = 1467
$ws.Cells.Item(56, 17).Value = 15

# Row 57
$ws.Cells.Item(57, 4).Value = 44811
$ws.Cells.Item(57, 9).Value = "Primera"
$ws.Cells.Item(57, 10).Value = 400
$ws.Cells.Item(57, 11).Value = 17000
$ws.Cells.Item(57, 12).Value = 18000
$ws.Cells.Item(57, 13).Value = 17425
$ws.Cells.Item(57, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(57, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(57, 16).Value = 1162
$ws.Cells.Item(57, 17).Value = 15

# Row 58
$ws.Cells.Item(58, 4).Value = 44392
$ws.Cells.Item(58, 9).Value = "Primera"
$ws.Cells.Item(58, 10).Value = 220
$ws.Cells.Item(58, 11).Value = 23000
$ws.Cells.Item(58, 12).Value = 23000
$ws.Cells.Item(58, 13).Value = 23000
$ws.Cells.Item(58, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(58, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(58, 16).Value = 1533
$ws.Cells.Item(58, 17).Value = 15

# Row 59
$ws.Cells.Item(59, 4).Value = 44406
$ws.Cells.Item(59, 9).Value = "Primera"
$ws.Cells.Item(59, 10).Value = 400
$ws.Cells.Item(59, 11).Value = 20000
$ws.Cells.Item(59, 12).Value = 22000
$ws.Cells.Item(59, 13).Value = 20850
$ws.Cells.Item(59, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(59, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(59, 16).Value = 1390
$ws.Cells.Item(59, 17).Value = 15

# Row 60
$ws.Cells.Item(60, 4).Value = 45043
$ws.Cells.Item(60, 9).Value = "Primera"
$ws.Cells.Item(60, 10).Value = 220
$ws.Cells.Item(60, 11).Value = 18000
$ws.Cells.Item(60, 12).Value = 20000
$ws.Cells.Item(60, 13).Value = 18909
$ws.Cells.Item(60, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(60, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(60, 16).Value = 1261
$ws.Cells.Item(60, 17).Value = 15

# Row 61
$ws.Cells.Item(61, 4).Value = 44776
$ws.Cells.Item(61, 9).Value = "Primera"
$ws.Cells.Item(61, 10).Value = 580
$ws.Cells.Item(61, 11).Value = 17000
$ws.Cells.Item(61, 12).Value = 19000
$ws.Cells.Item(61, 13).Value = 17897
$ws.Cells.Item(61, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(61, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(61, 16).Value = 1193
$ws.Cells.Item(61, 17).Value = 15

# Row 62
$ws.Cells.Item(62, 4).Value = 44832
$ws.Cells.Item(62, 9).Value = "Primera"
$ws.Cells.Item(62, 10).Value = 160
$ws.Cells.Item(62, 11).Value = 15000
$ws.Cells.Item(62, 12).Value = 15000
$ws.Cells.Item(62, 13).Value = 15000
$ws.Cells.Item(62, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(62, 15).Value = "Hijuelas"
$ws.Cells.Item(62, 16).Value = 1000
$ws.Cells.Item(62, 17).Value = 15

# Row 63
$ws.Cells.Item(63, 4).Value = 44832
$ws.Cells.Item(63, 9).Value = "Primera"
$ws.Cells.Item(63, 10).Value = 500
$ws.Cells.Item(63, 11).Value = 17000
$ws.Cells.Item(63, 12).Value = 19000
$ws.Cells.Item(63, 13).Value = 18080
$ws.Cells.Item(63, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(63, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(63, 16).Value = 1205
$ws.Cells.Item(63, 17).Value = 15

# Row 64
$ws.Cells.Item(64, 4).Value = 44742
$ws.Cells.Item(64, 9).Value = "Primera"
$ws.Cells.Item(64, 10).Value = 400
$ws.Cells.Item(64, 11).Value = 18000
$ws.Cells.Item(64, 12).Value = 20000
$ws.Cells.Item(64, 13).Value = 18850
$ws.Cells.Item(64, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(64, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(64, 16).Value = 1257
$ws.Cells.Item(64, 17).Value = 15

# Row 65
$ws.Cells.Item(65, 4).Value = 44817
$ws.Cells.Item(65, 9).Value = "Primera"
$ws.Cells.Item(65, 10).Value = 400
$ws.Cells.Item(65, 11).Value = 16000
$ws.Cells.Item(65, 12).Value = 17000
$ws.Cells.Item(65, 13).Value = 16425
$ws.Cells.Item(65, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(65, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(65, 16).Value = 1095
$ws.Cells.Item(65, 17).Value = 15

# Row 66
$ws.Cells.Item(66, 4).Value = 44817
$ws.Cells.Item(66, 9).Value = "Segunda"
$ws.Cells.Item(66, 10).Value = 150
$ws.Cells.Item(66, 11).Value = 15000
$ws.Cells.Item(66, 12).Value = 15000
$ws.Cells.Item(66, 13).Value = 15000
$ws.Cells.Item(66, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(66, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(66, 16).Value = 1000
$ws.Cells.Item(66, 17).Value = 15

# Row 67
$ws.Cells.Item(67, 4).Value = 45020
$ws.Cells.Item(67, 9).Value = "Primera"
$ws.Cells.Item(67, 10).Value = 240
$ws.Cells.Item(67, 11).Value = 22000
$ws.Cells.Item(67, 12).Value = 23000
$ws.Cells.Item(67, 13).Value = 22625
$ws.Cells.Item(67, 14).Value = "$/malla 17 kilos"
$ws.Cells.Item(67, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(67, 16).Value = 1331
$ws.Cells.Item(67, 17).Value = 17

# Row 68
$ws.Cells.Item(68, 4).Value = 44714
$ws.Cells.Item(68, 9).Value = "Primera"
$ws.Cells.Item(68, 10).Value = 200
$ws.Cells.Item(68, 11).Value = 16000
$ws.Cells.Item(68, 12).Value = 17000
$ws.Cells.Item(68, 13).Value = 16400
$ws.Cells.Item(68, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(68, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(68, 16).Value = 1093
$ws.Cells.Item(68, 17).Value = 15

# Row 69
$ws.Cells.Item(69, 4).Value = 44803
$ws.Cells.Item(69, 9).Value = "Primera"
$ws.Cells.Item(69, 10).Value = 400
$ws.Cells.Item(69, 11).Value = 16000
$ws.Cells.Item(69, 12).Value = 18000
$ws.Cells.Item(69, 13).Value = 16850
$ws.Cells.Item(69, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(69, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(69, 16).Value = 1123
$ws.Cells.Item(69, 17).Value = 15

# Row 70
$ws.Cells.Item(70, 4).Value = 44789
$ws.Cells.Item(70, 9).Value = "Primera"
$ws.Cells.Item(70, 10).Value = 400
$ws.Cells.Item(70, 11).Value = 15000
$ws.Cells.Item(70, 12).Value = 16000
$ws.Cells.Item(70, 13).Value = 15425
$ws.Cells.Item(70, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(70, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(70, 16).Value = 1028
$ws.Cells.Item(70, 17).Value = 15

# Row 71
$ws.Cells.Item(71, 4).Value = 44722
$ws.Cells.Item(71, 9).Value = "Primera"
$ws.Cells.Item(71, 10).Value = 150
$ws.Cells.Item(71, 11).Value = 18000
$ws.Cells.Item(71, 12).Value = 20000
$ws.Cells.Item(71, 13).Value = 18933
$ws.Cells.Item(71, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(71, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(71, 16).Value = 1262
$ws.Cells.Item(71, 17).Value = 15

# Row 72
$ws.Cells.Item(72, 4).Value = 44741
$ws.Cells.Item(72, 9).Value = "Primera"
$ws.Cells.Item(72, 10).Value = 250
$ws.Cells.Item(72, 11).Value = 18000
$ws.Cells.Item(72, 12).Value = 20000
$ws.Cells.Item(72, 13).Value = 18800
$ws.Cells.Item(72, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(72, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(72, 16).Value = 1253
$ws.Cells.Item(72, 17).Value = 15

# New row 73 (holds the data that used to belong to row 72)
$ws.Cells.Item(73, 1).Value = "6"
$ws.Cells.Item(73, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(73, 3).Value = "Metropolitana"
$ws.Cells.Item(73, 4).Value = 44398
$ws.Cells.Item(73, 5).Value = "13"
$ws.Cells.Item(73, 6).Value = "100112035"
$ws.Cells.Item(73, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(73, 8).Value = "Sin especificar"
$ws.Cells.Item(73, 9).Value = "Primera"
$ws.Cells.Item(73, 10).Value = 130
$ws.Cells.Item(73, 11).Value = 20000
$ws.Cells.Item(73, 12).Value = 20000
$ws.Cells.Item(73, 13).Value = 20000
$ws.Cells.Item(73, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(73, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(73, 16).Value = 1333
$ws.Cells.Item(73, 17).Value = 15
$ws.Cells.Item(73, 18).Value = "Hortaliza"

# Match the date-style formatting (s="2") used by the other rows in column D
$ws.Range("D73").NumberFormat = $ws.Range("D72").NumberFormat
